$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2: tweak a couple of account numbers
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A3").Value = "046FIA0015"
$ws2.Range("A5").Value = "046C000038"
$ws2.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# Sheet3: brand-new order-entry sheet, appended after Sheet2
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

$ws3.Columns.Item(1).ColumnWidth = 13.5703125
$ws3.Columns.Item(2).ColumnWidth = 20.85546875
$ws3.Columns.Item(3).ColumnWidth = 16.42578125
$ws3.Columns.Item(4).ColumnWidth = 13
$ws3.Columns.Item(5).ColumnWidth = 11.85546875
$ws3.Columns.Item(6).ColumnWidth = 10.42578125
$ws3.Columns.Item(7).ColumnWidth = 15.28515625

$ws3.Range("A1").Value = "af"
$ws3.Range("B1").Value = "orderType"
$ws3.Range("C1").Value = "symbol"
$ws3.Range("D1").Value = "qtty"
$ws3.Range("E1").Value = "price"
$ws3.Range("F1").Value = "BS"
$ws3.Range("G1").Value = "priceType"

$ws3.Range("A2").Value = "0001000474"
$ws3.Range("B2").Value = "Lệnh thông thường"
$ws3.Range("C2").Value = "AAA"
$ws3.Range("D2").Value = 100
$ws3.Range("E2").Value = 18
$ws3.Range("F2").Value = "Bán"
$ws3.Range("G2").Value = "LO"

# Header formatting: bold header style for A1:D1 (reuse sheet1/2 header look)
$ws3.Range("A1:D1").Font.Bold = $true
$ws3.Range("A1:D1").Interior.Pattern = -4124
$ws3.Range("A1:D1").Interior.ThemeColor = 10

# Header formatting for E1:G1: new font + vertical-centre alignment
$ws3.Range("E1:G1").Font.Name = "Roboto Mono Medium"
$ws3.Range("E1:G1").Font.Size = 9.8
$ws3.Range("E1:G1").Interior.Pattern = -4124
$ws3.Range("E1:G1").Interior.ThemeColor = 10
$ws3.Range("E1:G1").VerticalAlignment = -4108

$ws3.Range("B3").Select()
$ws3.Activate()

# ---------------------------------------------------------------------------
# Sheet1: update row 2's account-label cell, drop row 3 entirely
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C2").Value = "Margin. 046C000350 - Cá nhân trong nu?c 2 Cá nhân trong nu?c 2 Cá nhân trong nu?c 2"
$ws1.Rows.Item(3).Delete()
$ws1.Range("C11").Select()

$ws3.Activate()
